$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "5x2 Pinheader"
$ws.Range("D26").Value = "C358694"
$ws.Range("C26").Value = "TH"

$ws.Range("B26").Select()
